# Sort the data by F.Mov. (column A) ascending, then by Num. As. (column B) ascending.
# This matches the commit message: the sheet no longer sorts by Debe/Haber/Indice_Punteo,
# but instead by date (and assignment number as a tiebreaker).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$lastCol = 17  # column Q

$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1)), 0, 1, 0, 0) | Out-Null
$ws.Sort.SortFields.Add($ws.Range($ws.Cells.Item(2, 2), $ws.Cells.Item($lastRow, 2)), 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()
